$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04172184405617529
$ws.Range("C2").Value = 0.04103571897497393
$ws.Range("D2").Value = 3.223369029078222
$ws.Range("E2").Value = 13.86384647080068
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 17.16997306291006
